$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the existing "Values" comment (currently anchored at A2 / the "red" row) ---
$valuesComment = $ws.Range("A2").Comment
$valuesCommentText = $valuesComment.Text()
$valuesComment.Delete()

# --- row 2 becomes the new "test_color" row ---
$ws.Range("A2").Value = "test_color"
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("B2").Value = 180
$ws.Range("C2").Value = 50

# --- header row: add new "type" column (after test_color so shared-string order matches) ---
$ws.Range("D1").Value = "type"
$ws.Range("D2").Value = "low"

# --- row 3 stays "green", gains a "high" type ---
$ws.Range("A3").Value = "green"
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("B3").Value = 120
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "high"

# --- row 4 (new) becomes the "red" row, also "high" type ---
$ws.Range("A4").Value = "red"
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = "high"
$ws.Rows(4).RowHeight = 20.1

# --- re-anchor the "Values" comment on the row that now holds "red" ---
$newComment = $ws.Range("A4").AddComment($valuesCommentText)

# --- update selection to follow the new last row ---
$ws.Range("D5").Select()
